$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells hold text that looks numeric (e.g. "241.85"); Excel
# auto-converts plain .Value assignment of such text to a Number. To
# preserve the original text semantics we briefly mark the cell as
# Text ("@"), assign the literal string, then restore the Normal style
# so no stray number-format style is left attached to the cell.
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '29.320.49'
$ws.Range("E2").Value = '  +0.37%  '

Set-TextValue $ws.Range("D3") '1.875.39'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("E4").Value = '  +0.03%  '

Set-TextValue $ws.Range("D5") '0.7130'

Set-TextValue $ws.Range("D6") '241.85'
$ws.Range("E6").Value = '  +0.61%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +0.95%  '

Set-TextValue $ws.Range("D9") '0.07747'
$ws.Range("E9").Value = '  -0.18%  '

Set-TextValue $ws.Range("D10") '24.90'
$ws.Range("E10").Value = '  -0.85%  '

Set-TextValue $ws.Range("D11") '0.08509'
$ws.Range("E11").Value = '  +2.94%  '

Set-TextValue $ws.Range("D12") '1.881.37'
$ws.Range("E12").Value = '  +1.08%  '

Set-TextValue $ws.Range("D13") '5.216'
$ws.Range("E13").Value = '  -0.23%  '

Set-TextValue $ws.Range("D14") '0.7099'
$ws.Range("E14").Value = '  -0.87%  '

Set-TextValue $ws.Range("D15") '91.50'
$ws.Range("E15").Value = '  +1.24%  '

Set-TextValue $ws.Range("D16") '29.313.21'
$ws.Range("E16").Value = '  +0.42%  '

Set-TextValue $ws.Range("D17") '0.000008194'
$ws.Range("E17").Value = '  +5.12%  '

Set-TextValue $ws.Range("D18") '6.005'
$ws.Range("E18").Value = '  +2.47%  '

Set-TextValue $ws.Range("D19") '241.82'
$ws.Range("E19").Value = '  -0.93%  '

Set-TextValue $ws.Range("D20") '2.134.63'
$ws.Range("E20").Value = '  +1.62%  '

$ws.Range("E21").Value = '  +0.80%  '

Set-TextValue $ws.Range("D22") '0.9997'

Set-TextValue $ws.Range("D23") '7.811'
$ws.Range("E23").Value = '  -1.94%  '

Set-TextValue $ws.Range("D24") '0.9999'
$ws.Range("E24").Value = '  -0.01%  '

Set-TextValue $ws.Range("D25") '0.1605'
$ws.Range("E25").Value = '  +0.92%  '

Set-TextValue $ws.Range("D26") '162.92'
$ws.Range("E26").Value = '  +0.23%  '

Set-TextValue $ws.Range("D27") '9.050'
$ws.Range("E27").Value = '  +1.61%  '

Set-TextValue $ws.Range("D28") '18.48'
$ws.Range("E28").Value = '  +0.81%  '

Set-TextValue $ws.Range("D29") '1.513'
$ws.Range("E29").Value = '  +1.16%  '

Set-TextValue $ws.Range("D30") '4.397'
$ws.Range("E30").Value = '  -0.45%  '

Set-TextValue $ws.Range("D31") '4.318'
$ws.Range("E31").Value = '  +2.12%  '

Set-TextValue $ws.Range("D32") '1.284'
$ws.Range("E32").Value = '  -2.35%  '

Set-TextValue $ws.Range("D33") '0.05258'
$ws.Range("E33").Value = '  +1.37%  '

Set-TextValue $ws.Range("D34") '1.932'
$ws.Range("E34").Value = '  +1.14%  '

$ws.Range("E35").Value = '  +0.27%  '

Set-TextValue $ws.Range("D36") '0.7422'
$ws.Range("E36").Value = '  +2.18%  '

Set-TextValue $ws.Range("D37") '2.687'
$ws.Range("E37").Value = '  +0.49%  '

Set-TextValue $ws.Range("D38") '0.01864'
$ws.Range("E38").Value = '  +0.53%  '

Set-TextValue $ws.Range("D39") '2.719'

Set-TextValue $ws.Range("D40") '1.184.57'
$ws.Range("E40").Value = '  +1.71%  '

Set-TextValue $ws.Range("D41") '6.380'
$ws.Range("E41").Value = '  +3.61%  '

Set-TextValue $ws.Range("D42") '0.8878'
$ws.Range("E42").Value = '  -1.86%  '

Set-TextValue $ws.Range("D43") '72.87'
$ws.Range("E43").Value = '  +0.94%  '

Set-TextValue $ws.Range("D44") '106.24'
$ws.Range("E44").Value = '  +4.57%  '

$ws.Range("E45").Value = '  -0.05%  '

Set-TextValue $ws.Range("D46") '2.030.39'
$ws.Range("E46").Value = '  +1.49%  '

$ws.Range("E47").Value = '  +2.55%  '

$ws.Range("E48").Value = '  -0.20%  '

$ws.Range("E49").Value = '  +1.24%  '

Set-TextValue $ws.Range("D50") '9.375'
$ws.Range("E50").Value = '  +0.64%  '

Set-TextValue $ws.Range("D51") '0.4312'
$ws.Range("E51").Value = '  +1.14%  '
